$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
$cs = $nm.Theme.ThemeColorScheme
$cs.Item(4).RGB = 11111111
